# Edit: update the roadmap date on the title slide, and remove the
# "Results: C++ to C# via staged YAML" slide from the deck.

$p = $ppt.ActivePresentation

# 1) Update the date shown on the title slide from "11-19-2023" to
#    "11-29-2023", leaving everything else about the run/paragraph intact.
$oldDate = "11-19-2023"
$newDate = "11-29-2023"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            $fullText = $tr.Text
            $idx = $fullText.IndexOf($oldDate)
            if ($idx -ge 0) {
                $dateRange = $tr.Characters($idx + 1, $oldDate.Length)
                $dateRange.Text = $newDate
            }
        }
    }
}

# 2) Delete the "Results: C++ to C# via staged YAML" slide from the deck.
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $slide = $p.Slides.Item($i)
    $isTarget = $false
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "Results: C++ to C# via staged YAML") {
                $isTarget = $true
                break
            }
        }
    }
    if ($isTarget) {
        $slide.Delete()
    }
}
